# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp string (cell A1)
$ws.Range("A1").Value = "Datos actualizados a 2 de Junio de 2020 a las 00:35"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1857586
$ws.Range("C4").Value = 20416
$ws.Range("D4").Value = 608487
$ws.Range("E4").Value = 1142221
$ws.Range("G4").Value = 683
$ws.Range("H4").Value = 106878

# Row 5 - Brasil
$ws.Range("B5").Value = 525307
$ws.Range("C5").Value = 10458
$ws.Range("E5").Value = 288975
$ws.Range("G5").Value = 463
$ws.Range("H5").Value = 29777

# Row 34 - Colombia
$ws.Range("B34").Value = 30493
$ws.Range("C34").Value = 1110
$ws.Range("D34").Value = 9661
$ws.Range("E34").Value = 19863
$ws.Range("G34").Value = 30
$ws.Range("H34").Value = 969

# Row 61 - Noruega
$ws.Range("B61").Value = 8446
$ws.Range("C61").Value = 6
$ws.Range("E61").Value = 483

# Row 129 - Principado de Andorra
$ws.Range("B129").Value = 765
$ws.Range("C129").Value = 1
$ws.Range("D129").Value = 698
$ws.Range("E129").Value = 16
